$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Drop all existing hyperlinks up front. In this engine, deleting
#    a hyperlink collection scoped to any single range actually clears
#    every hyperlink on the sheet, so do it exactly once before the
#    rows move around; we'll re-add the 8 links at their final
#    (post-insert) addresses at the end.
# ------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()

# ------------------------------------------------------------------
# 2. Insert one blank row after each "<Pose>VideoPath" row (just before
#    the matching "<Pose>VideoDescription" row) so a new
#    "<Pose>OnlineID" row can be written into it. Working from the
#    bottom-most insertion point upward means every Insert() call uses
#    the same row number it would have had in the original sheet,
#    since nothing below has shifted yet when we act on it.
# ------------------------------------------------------------------
$ws.Rows.Item(24).Insert()
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(20).Insert()
$ws.Rows.Item(18).Insert()
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(12).Insert()

# ------------------------------------------------------------------
# 3. Populate the freshly inserted rows with the YouTube-ID label/value
#    pairs, matching the existing "TreePoseOnlineID" row's formatting
#    (plain label cell in column A, wrap-text value cell in column B).
# ------------------------------------------------------------------
function Set-OnlineIdRow($row, $label, $value) {
    $ws.Range("A$row").Value = $label
    $ws.Range("B$row").Value = $value
    $ws.Range("B$row").WrapText = $true
}

Set-OnlineIdRow 12 "CatPoseOnlineID" "d8M6uOmp_9o"
Set-OnlineIdRow 15 "BowPoseOnlineID" "NMSGa2NuIzc"
Set-OnlineIdRow 18 "BridgePoseOnlineID" "zM--m3JOSSs"
Set-OnlineIdRow 21 "CamelPoseOnlineID" "nsT9naztI_I"
Set-OnlineIdRow 24 "ChairPoseOnlineID" "2Xo0PQHlMh0"
Set-OnlineIdRow 27 "CowPoseOnlineID" "IQsgImpibUo"
Set-OnlineIdRow 30 "TrianglePoseOnlineID" "B169KXlIFak"

# ------------------------------------------------------------------
# 4. Re-create the hyperlinks at their shifted addresses (video-path
#    rows move down by 1 for every "OnlineID" row inserted above them).
#    Hyperlinks.Add() always stamps a fresh "applied hyperlink font"
#    style onto its target cell, so immediately reset each cell back
#    to the plain named "Hyperlink" style to match the original
#    formatting instead of layering a duplicate on top of it.
# ------------------------------------------------------------------
function Add-PoseHyperlink($row, $fileName) {
    $ws.Hyperlinks.Add($ws.Range("B$row"), "file:///\\Videos\$fileName")
    $ws.Range("B$row").Style = "Hyperlink"
}

Add-PoseHyperlink 5 "TreePose.mp4"
Add-PoseHyperlink 11 "CatPose.mp4"
Add-PoseHyperlink 14 "BowPose.mp4"
Add-PoseHyperlink 17 "BridgePose.mp4"
Add-PoseHyperlink 20 "CamelPose.mp4"
Add-PoseHyperlink 23 "ChairPose.mp4"
Add-PoseHyperlink 26 "CowPose.mp4"
Add-PoseHyperlink 29 "TrianglePose.mp4"

# ------------------------------------------------------------------
# 5. Leave the selection on the last edited cell, like the author did.
# ------------------------------------------------------------------
$ws.Range("B30").Select()
